# Update "想去人数" (attendee count) figures for four events.
# These values live on both the "展览" sheet and the duplicate "全部类型"
# sheet, so apply the same updates to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 357
    $ws.Range("F4").Value = 1542
    $ws.Range("F8").Value = 58
    $ws.Range("F9").Value = 383
}
